# Preparing for the next public release.
# Update the "unidentified" plankton-group labels to the consolidated
# " Other microalgae" group, and add three new "phylum incertae sedis"
# taxa rows under the Phylum section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rows 2 & 3: collapse the bespoke "Unidentified ... eukaryots" plankton
#    group labels into the shared " Other microalgae" group.
$ws.Range("C2").Value = " Other microalgae"
$ws.Range("C3").Value = " Other microalgae"

# 2) Insert 4 rows before the old row 14 (Dinophyceae/Class...) so that the
#    "Phylum" block gains 3 new taxa rows (13-15) and keeps its usual blank
#    separator row before the "Class" block (now rows 16-17).
$ws.Rows("13:16").Insert()

$ws.Range("A13").Value = "Flagellates phylum incertae sedis"
$ws.Range("B13").Value = "Phylum"
$ws.Range("C13").Value = " Other microalgae"

$ws.Range("A14").Value = "Eukarotic picoplankton phylum incertae sedis"
$ws.Range("B14").Value = "Phylum"
$ws.Range("C14").Value = " Other microalgae"

$ws.Range("A15").Value = "Unicells phylum incertae sedis"
$ws.Range("B15").Value = "Phylum"
$ws.Range("C15").Value = " Other microalgae"

# 3) Restore the UI selection marker left by the author at save time.
$ws.Range("C43").Select()
